$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.73810887336731
$ws.Range("B1").Value = 3.701992750167847
$ws.Range("C1").Value = 2.27217435836792
$ws.Range("D1").Value = 1.595670700073242
$ws.Range("E1").Value = 1.358129739761353
